$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "18×11=198" "22×61=1342"
Replace-Text "13×14=182" "71×75=5325"
Replace-Text "48×91=4368" "76×99=7524"
Replace-Text "67×11=737" "74×74=5476"
Replace-Text "79×67=5293" "68×77=5236"
Replace-Text "24×26=624" "56×18=1008"
Replace-Text "47×56=2632" "84×26=2184"
Replace-Text "54×55=2970" "65×45=2925"
Replace-Text "16×86=1376" "78×17=1326"
Replace-Text "48×42=2016" "98×81=7938"
Replace-Text "44×77=3388" "77×92=7084"
Replace-Text "76×63=4788" "14×21=294"
Replace-Text "93×56=5208" "82×57=4674"
Replace-Text "86×33=2838" "82×89=7298"
Replace-Text "99×81=8019" "32×41=1312"
Replace-Text "33×23=759" "96×13=1248"
Replace-Text "77×90=6930" "61×59=3599"
Replace-Text "94×31=2914" "83×81=6723"
Replace-Text "86×77=6622" "26×67=1742"
Replace-Text "84×37=3108" "87×26=2262"
Replace-Text "87×63=5481" "85×26=2210"
Replace-Text "96×61=5856" "31×51=1581"
Replace-Text "79×87=6873" "36×91=3276"
Replace-Text "76×91=6916" "81×70=5670"
Replace-Text "60×53=3180" "45×71=3195"
